# Add files via upload
#
# The author inserted a new row above the data table (pushing the header
# row and all data rows down by one) and used that new row to record a
# note about an outlier that was removed from the SY1140A readings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 7 (everything from the old row 7 onward
# shifts down by one row, carrying its formatting/formulas with it).
$ws.Rows.Item(7).Insert()

# Record the outlier note in the newly-inserted row.
$ws.Range("A7").Value = "Outliers were removed prior to data analysis. One outliers from SY1140A."

# Leave the freshly-edited cell selected, matching the saved view state.
$ws.Range("A7").Select()
